# Updated cryptos list on Thu Oct 31 11:53:42 UTC 2024 with GitHub Actions
#
# Refreshes the price (column D) and 1h-volume-change (column E) figures for
# the coinranking.com snapshot on Sheet1, and fixes the Uniswap/BitcoinCash
# rows which had swapped rank positions (row 20 <-> row 21).
#
# All writes use a leading apostrophe so Excel stores the figures as literal
# text (matching how the sheet already represents prices such as "1.00" or
# "72.267.49") instead of silently coercing them to numbers and dropping
# formatting such as trailing zeros or the thousands-dot grouping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'72.267.49"
$ws.Range("E2").Value = "'  +0.27%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'2.637.51"
$ws.Range("E3").Value = "'  -0.94%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "'  -0.05%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'586.17"
$ws.Range("E5").Value = "'  -1.85%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'175.06"
$ws.Range("E6").Value = "'  -0.54%  "

# Row 7 - USDC
$ws.Range("E7").Value = "'  -0.06%  "

# Row 8 - XRP
$ws.Range("E8").Value = "'  -0.44%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "'2.635.52"
$ws.Range("E9").Value = "'  -1.00%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "'  +1.65%  "

# Row 11 - TRON
$ws.Range("E11").Value = "'  +1.45%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "'  +1.93%  "

# Row 13 - Toncoin
$ws.Range("E13").Value = "'  -1.27%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "'3.118.99"
$ws.Range("E14").Value = "'  -1.06%  "

# Row 15 - ShibaInu
$ws.Range("E15").Value = "'  +0.42%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "'72.146.19"
$ws.Range("E16").Value = "'  +0.22%  "

# Row 17 - Avalanche
$ws.Range("D17").Value = "'25.75"
$ws.Range("E17").Value = "'  -1.64%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "'2.622.60"
$ws.Range("E18").Value = "'  -1.62%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "'12.08"
$ws.Range("E19").Value = "'  +0.69%  "

# Row 20 - was BitcoinCash, is now Uniswap (rows 20/21 swapped rank)
$ws.Range("B20").Value = "'Uniswap"
$ws.Range("C20").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'7.97"
$ws.Range("E20").Value = "'  -0.32%  "

# Row 21 - was Uniswap, is now BitcoinCash
$ws.Range("B21").Value = "'BitcoinCash"
$ws.Range("C21").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'375.16"
$ws.Range("E21").Value = "'  +1.34%  "

# Row 22 - Polkadot
$ws.Range("E22").Value = "'  -1.21%  "

# Row 23 - SuiNetwork
$ws.Range("E23").Value = "'  +0.08%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'71.46"
$ws.Range("E24").Value = "'  -0.18%  "

# Row 25 - Dai
$ws.Range("E25").Value = "'  -0.05%  "

# Row 26 - NEARProtocol
$ws.Range("D26").Value = "'4.22"
$ws.Range("E26").Value = "'  -1.93%  "

# Row 27 - Aptos
$ws.Range("D27").Value = "'9.41"
$ws.Range("E27").Value = "'  -3.82%  "

# Row 28 - WrappedeETH
$ws.Range("D28").Value = "'2.771.83"
$ws.Range("E28").Value = "'  -1.09%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "'  -0.02%  "

# Row 30 - PEPE (price uses a subscript-3 in the leading-zero count)
$ws.Range("D30").Value = "'0.0" + [char]0x2083 + "0947"
$ws.Range("E30").Value = "'  +1.45%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("E31").Value = "'  -1.17%  "

# Row 32 - Bittensor
$ws.Range("D32").Value = "'490.85"
$ws.Range("E32").Value = "'  -3.33%  "

# Row 33 - Fetch.AI
$ws.Range("D33").Value = "'1.30"
$ws.Range("E33").Value = "'  +1.43%  "

# Row 34 - PancakeSwap
$ws.Range("D34").Value = "'1.80"
$ws.Range("E34").Value = "'  -0.62%  "

# Row 35 - FirstDigitalUSD
$ws.Range("E35").Value = "'  -0.04%  "

# Row 36 - Monero
$ws.Range("D36").Value = "'160.09"
$ws.Range("E36").Value = "'  -3.00%  "

# Row 37 - Kaspa
$ws.Range("D37").Value = "'0.115"
$ws.Range("E37").Value = "'  +7.89%  "

# Row 38 - EthereumClassic
$ws.Range("D38").Value = "'19.16"
$ws.Range("E38").Value = "'  -1.78%  "

# Row 40 - ImmutableX
$ws.Range("E40").Value = "'  -0.90%  "

# Row 41 - USDe
$ws.Range("E41").Value = "'  +0.00%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "'  -4.04%  "

# Row 43 - dogwifhat
$ws.Range("D43").Value = "'2.58"
$ws.Range("E43").Value = "'  +1.20%  "

# Row 44 - RenderToken
$ws.Range("E44").Value = "'  -2.23%  "

# Row 45 - PolygonEcosystemToken
$ws.Range("E45").Value = "'  -1.71%  "

# Row 46 - OKB
$ws.Range("D46").Value = "'39.01"
$ws.Range("E46").Value = "'  -0.54%  "

# Row 47 - Aave
$ws.Range("D47").Value = "'150.18"
$ws.Range("E47").Value = "'  -1.49%  "

# Row 48 - Filecoin
$ws.Range("E48").Value = "'  -2.28%  "

# Row 49 - ARBITRUM
$ws.Range("E49").Value = "'  -0.40%  "

# Row 50 - Optimism
$ws.Range("E50").Value = "'  -2.59%  "

# Row 51 - Mantle
$ws.Range("D51").Value = "'0.606"
$ws.Range("E51").Value = "'  +1.09%  "
